# Update "想去人数" (want-to-go count) values in sheet "展览" (Exhibition)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 67
$ws1.Range("F3").Value = 656
$ws1.Range("F4").Value = 230
$ws1.Range("F6").Value = 9919
$ws1.Range("F7").Value = 898
$ws1.Range("F9").Value = 1239
$ws1.Range("F10").Value = 4458
$ws1.Range("F11").Value = 2
$ws1.Range("F15").Value = 57
$ws1.Range("F17").Value = 290
$ws1.Range("F18").Value = 566
$ws1.Range("F19").Value = 110
$ws1.Range("F21").Value = 1481

# Update sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 20

# Update sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 67
$ws4.Range("F3").Value = 20
$ws4.Range("F4").Value = 656
$ws4.Range("F5").Value = 230
$ws4.Range("F7").Value = 9919
$ws4.Range("F8").Value = 898
$ws4.Range("F10").Value = 1239
$ws4.Range("F11").Value = 4458
$ws4.Range("F12").Value = 2
$ws4.Range("F16").Value = 57
$ws4.Range("F18").Value = 290
$ws4.Range("F19").Value = 566
$ws4.Range("F20").Value = 110
$ws4.Range("F22").Value = 1481
